$d = $word.ActiveDocument

# 1. Fix "BI" -> "BL" in the heading text.
$d.Content.Find.Execute("Requests to BI layer:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Requests to BL layer:", 2)

# 2. Move the "_GoBack" bookmark so that it sits right after "Requests to BL"
#    (before the " layer:" text) instead of its old spot inside "SaveUsersData".
#    Re-adding a bookmark with the same name automatically removes any
#    previous bookmark of that name elsewhere in the document.
$p1 = $d.Paragraphs(1)
$headingStart = $p1.Range.Start
$markerPos = $headingStart + "Requests to BL".Length
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
